# Project Euler 006 workbook - presentation fix
#
# The original layout used columns C:E (data in C:D, "<-- answer" label in
# E2) starting at row 2, with no title. The new layout adds a title row and
# shifts everything down two rows and left one column (B:D, data starting
# row 7, header row 6, the SUM()^2-SUM() result / "<-- answer" on row 4, and
# a new title in B2).
#
# Achieve the row/column shift the same way a human would in Excel: insert
# two blank rows at the top, then delete the (now blank) original first
# data column so everything slides one column to the left. Formulas are
# auto-adjusted by Excel's insert/delete, so SUM(C5:C104) etc. become
# SUM(B7:B106) etc. automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Push everything down by 2 rows (new blank rows 1:2 inserted at top).
$ws.Rows("1:2").Insert()

# 2. Push everything left by 1 column (delete the now-empty former column B).
$ws.Columns("B:B").Delete()

# 3. Add the new title in B2, using the same number format / bold style as
#    the neighbouring "<-- answer" result cell (style index 2 in the
#    original file) so no redundant style gets created.
$ws.Range("B2").Value = "Project Euler 6: Sum Square Difference"
$ws.Range("C4").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 4. Match the saved selection state (whole column B selected, as if the
#    user had just clicked the column header after typing the title).
$ws.Columns("B:B").Select()
